$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6, shifting existing rows 6..38 down to 7..39
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new data record
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(6, 3).Value = "Los Lagos"
$ws.Cells.Item(6, 4).Value = 44462
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = 100112026
$ws.Cells.Item(6, 7).Value = "Haba"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(6, 11).Value = 15000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 13).Value = 15000
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 16).Value = 600
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Apply the same number format as the D column date cells to the new D6 cell
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
